$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking values
# (e.g. "1.00", "28.80", "493.50") keep their exact original formatting
# instead of being coerced into Excel numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Price (column D) updates ---
$ws.Range("D2").Value = "76.342.12"
$ws.Range("D3").Value = "3.043.77"
$ws.Range("D5").Value = "198.45"
$ws.Range("D6").Value = "617.19"
$ws.Range("D10").Value = "3.044.01"
$ws.Range("D11").Value = "0.438"
$ws.Range("D12").Value = "0.161"
$ws.Range("D14").Value = "3.599.85"
$ws.Range("D15").Value = "28.80"
$ws.Range("D16").Value = "76.253.66"
$ws.Range("D18").Value = "3.038.61"
$ws.Range("D19").Value = "13.51"
$ws.Range("D21").Value = "382.57"
$ws.Range("D22").Value = "2.37"
$ws.Range("D24").Value = "3.195.70"
$ws.Range("D26").Value = "0.999"
$ws.Range("D28").Value = "9.73"
$ws.Range("D30").Value = "1.00"
$ws.Range("D33").Value = "493.50"
$ws.Range("D35").Value = "1.00"
$ws.Range("D36").Value = "20.57"
$ws.Range("D37").Value = "162.23"
$ws.Range("D38").Value = "0.118"
$ws.Range("D40").Value = "191.64"
$ws.Range("D41").Value = "0.380"
$ws.Range("D44").Value = "0.792"
$ws.Range("D46").Value = "41.97"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("E3").Value = "  +3.74%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("E6").Value = "  +3.38%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("E9").Value = "  +4.04%  "
$ws.Range("E10").Value = "  +3.89%  "
$ws.Range("E11").Value = "  -2.32%  "
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("E13").Value = "  +5.91%  "
$ws.Range("E14").Value = "  +3.94%  "
$ws.Range("E15").Value = "  +2.81%  "
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("E18").Value = "  +4.17%  "
$ws.Range("E19").Value = "  +1.66%  "
$ws.Range("E20").Value = "  +2.28%  "
$ws.Range("E21").Value = "  +2.72%  "
$ws.Range("E22").Value = "  +2.92%  "
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("E24").Value = "  +4.51%  "
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  +1.50%  "
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("E31").Value = "  +4.65%  "
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("E34").Value = "  +4.43%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("E38").Value = "  +6.82%  "
$ws.Range("E39").Value = "  +1.85%  "
$ws.Range("E40").Value = "  +7.15%  "
$ws.Range("E41").Value = "  -2.27%  "
$ws.Range("E42").Value = "  -4.40%  "
$ws.Range("E44").Value = "  +20.32%  "
$ws.Range("E45").Value = "  +3.45%  "
$ws.Range("E47").Value = "  +4.28%  "
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("E49").Value = "  +4.80%  "
$ws.Range("E50").Value = "  +2.53%  "
$ws.Range("E51").Value = "  -0.25%  "
